$d = $word.ActiveDocument

$replacements = @(
    @{Old = "306÷4="; New = "985÷9="},
    @{Old = "743÷6="; New = "663÷7="},
    @{Old = "360÷4="; New = "835÷8="},
    @{Old = "923÷2="; New = "424÷2="},
    @{Old = "473÷2="; New = "333÷4="},
    @{Old = "984÷4="; New = "647÷8="},
    @{Old = "278÷9="; New = "493÷2="},
    @{Old = "230÷3="; New = "771÷6="},
    @{Old = "429÷9="; New = "548÷8="},
    @{Old = "604÷2="; New = "781÷3="},
    @{Old = "335÷2="; New = "769÷2="},
    @{Old = "728÷3="; New = "870÷5="},
    @{Old = "197÷7="; New = "419÷5="},
    @{Old = "562÷9="; New = "387÷2="},
    @{Old = "124÷4="; New = "213÷7="},
    @{Old = "151÷4="; New = "807÷5="},
    @{Old = "301÷9="; New = "529÷4="},
    @{Old = "576÷3="; New = "779÷6="},
    @{Old = "817÷5="; New = "335÷7="},
    @{Old = "857÷6="; New = "925÷8="},
    @{Old = "948÷2="; New = "447÷8="},
    @{Old = "676÷4="; New = "584÷7="},
    @{Old = "556÷4="; New = "843÷9="},
    @{Old = "458÷6="; New = "730÷6="},
    @{Old = "593÷3="; New = "207÷2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
